$wb = $excel.ActiveWorkbook

# --- Content edits -------------------------------------------------------
# "verifySearchResultsOnTyping": A1 "TypeJointUse" -> "TypePoleReplacement",
#                                  A2 "Joint Use" -> "Pole Replacement"
# (written first so new shared-string entries are interned in the same
#  order the target workbook uses)
$wsTyping = $wb.Worksheets.Item("verifySearchResultsOnTyping")
$wsTyping.Range("A1").Value = "TypePoleReplacement"
$wsTyping.Range("A2").Value = "Pole Replacement"

# "verifyDynamicDropdown": A2 "Ind" -> "Indus"
$wsDropdown = $wb.Worksheets.Item("verifyDynamicDropdown")
$wsDropdown.Range("A2").Value = "Indus"

# --- Selection / active sheet edits --------------------------------------
# The selection on "verifySearchResultsOnTyping" moves from A1:A2 to just A2.
$wsTyping.Range("A2").Select()

# The active/selected tab moves from "verifyExternalVideo" (last sheet) to
# "verifyDynamicDropdown".
$wsDropdown.Activate()

Write-Output "edits applied"
